$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the "Testing:" and "Demo:" end dates out further, which recalculates
# the "Time Spent (days)" formulas in column E (and the chart that reads it).
$ws.Range("D10").Value = "5/17/2020"
$ws.Range("D11").Value = "5/17/2020"

# Column D needs to grow to fit the new, wider dates.
$ws.Columns.Item(4).AutoFit() | Out-Null

# Update the view state: scroll down one row and move the selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("I10").Select() | Out-Null
